$wb = $excel.ActiveWorkbook

$ws5 = $wb.Worksheets.Item("Sheet5")
$ws6 = $wb.Worksheets.Item("Sheet6")

# Set new string values in the order they were first introduced so the
# shared string table is rebuilt with the same ordering as the target:
#   10 = "q7", 11 = "l", 12 = "r"
$ws5.Range("A4").Value = "q7"
$ws6.Range("A6").Value = "l"
$ws5.Range("A6").Value = "r"

# Remaining cells (reuse existing shared strings / plain numbers)
$ws5.Range("B4").Value = 0.0
$ws5.Range("A5").Value = "a"
$ws5.Range("B5").Value = 1563.0
$ws5.Range("B6").Value = 2084.0

$ws6.Range("A5").Value = "a"
$ws6.Range("B5").Value = 4168.0
$ws6.Range("B6").Value = 0.0
